$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 227, shifting existing rows 227-274 down to 228-275.
$ws.Rows("227").Insert()

# Populate the newly inserted row 227 with the new record.
$ws.Cells.Item(227, 1).Value = 1
$ws.Cells.Item(227, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(227, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(227, 4).Value = 44782
$ws.Cells.Item(227, 5).Value = 15
$ws.Cells.Item(227, 6).Value = "Fruta"
$ws.Cells.Item(227, 7).Value = 100108
$ws.Cells.Item(227, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(227, 9).Value = 100108006
$ws.Cells.Item(227, 10).Value = "Plátano"
$ws.Cells.Item(227, 11).Value = "Sin especificar"
$ws.Cells.Item(227, 12).Value = "Pintón"
$ws.Cells.Item(227, 13).Value = 120
$ws.Cells.Item(227, 14).Value = 22000
$ws.Cells.Item(227, 15).Value = 23000
$ws.Cells.Item(227, 16).Value = 22500
$ws.Cells.Item(227, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(227, 18).Value = "Ecuador"
$ws.Cells.Item(227, 19).Value = 1125
$ws.Cells.Item(227, 20).Value = 20
